# ==========================================================================
# Adds a new "2022-Q4" quarter: a summary row in "总计" and a brand-new
# "2022-Q4" worksheet (holding positions) inserted right after "总计" and
# before "2022-Q3". All the other quarter sheets shift right but keep their
# own content unchanged.
# ==========================================================================

$wb = $excel.ActiveWorkbook

# --- helper: make a range look like the bold/bordered/centered header-ish
#     style ("s=2" in the original sheets) used for the header row and the
#     index column. -------------------------------------------------------
function Set-IndexStyle($rng) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108   # xlCenter
    $rng.VerticalAlignment = -4160     # xlTop
    $rng.Borders.LineStyle = 1         # xlContinuous
}

# --- helper: assign a numeric-looking string as TEXT (keeps leading zeros,
#     avoids Excel's automatic "looks like a number" coercion). -----------
function Set-TextValue($rng, $text) {
    $rng.Value = "'" + $text
}

# ---------------------------------------------------------------------
# 1. "总计" sheet: insert a new row 2 for 2022-Q4, push the rest down.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").ClearFormats()

$wsTotal.Range("A2").Value = 0
Set-IndexStyle $wsTotal.Range("A2")
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 4
$wsTotal.Range("D2").Value = 0.2

# Re-sequence the A-index column (0-based row counter) for the rows that
# got pushed down, since their old index values came along for the ride.
for ($r = 3; $r -le 7; $r++) {
    $wsTotal.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2. Brand-new "2022-Q4" worksheet, positioned right after "总计".
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q3")
$wsQ4 = $wb.Worksheets.Add($beforeSheet)
$wsQ4.Name = "2022-Q4"

# Header row.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $wsQ4.Cells.Item(1, $i + 2)   # starts at column B
    $cell.Value = $headers[$i]
}
Set-IndexStyle $wsQ4.Range("B1:H1")

# Data rows: A, B(code,text), C(name,text), D(text), E(text), F(text), G(text), H(number)
$rows = @(
    @(0, "160218", "国泰国证房地产行业指数A",       "4.76", "92.61", "2.14", "0.1019", 8),
    @(1, "515760", "华夏中证浙江国资创新发展ETF",   "2.14", "99.05", "3.62", "0.0775", 9),
    @(2, "015042", "国泰国证房地产行业指数C",       "0.79", "92.61", "2.14", "0.0169", 8),
    @(3, "510190", "华安上证龙头ETF",               "0.58", "98.19", "1.13", "0.0066", 2)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $wsQ4.Cells.Item($r, 1).Value = $row[0]
    Set-IndexStyle $wsQ4.Cells.Item($r, 1)

    Set-TextValue $wsQ4.Cells.Item($r, 2) $row[1]
    $wsQ4.Cells.Item($r, 3).Value = $row[2]
    Set-TextValue $wsQ4.Cells.Item($r, 4) $row[3]
    Set-TextValue $wsQ4.Cells.Item($r, 5) $row[4]
    Set-TextValue $wsQ4.Cells.Item($r, 6) $row[5]
    Set-TextValue $wsQ4.Cells.Item($r, 7) $row[6]
    $wsQ4.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------
# 3. Keep the originally-active sheet selected (last tab, "2021-Q2").
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()

Write-Host "2022-Q4 sheet added; workbook now has $($wb.Worksheets.Count) sheets."
